$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.414658091557881
$ws.Range("D2").Value = 3.710838819035526
$ws.Range("E2").Value = 16.58781092826514
$ws.Range("F2").Value = 19.07606382731735
$ws.Range("G2").Value = 20.08707642018283
$ws.Range("H2").Value = 11.87795694170807
$ws.Range("I2").Value = 17.79008391292007
$ws.Range("K2").Value = 11.30786179471816
$ws.Range("N2").Value = 16.75202072264581
$ws.Range("O2").Value = 16.93833874414826
$ws.Range("B3").Value = 6.289114304817435
$ws.Range("D3").Value = 3.636313135138622
$ws.Range("E3").Value = 15.64205953995149
$ws.Range("F3").Value = 19.02453921908501
$ws.Range("G3").Value = 19.95151561870688
$ws.Range("H3").Value = 11.91086151460003
$ws.Range("I3").Value = 17.88501248978983
$ws.Range("K3").Value = 10.86349663151118
$ws.Range("N3").Value = 16.79372421800811
$ws.Range("O3").Value = 16.96298881989019
$ws.Range("B4").Value = 6.211643869133904
$ws.Range("D4").Value = 3.589159604468868
$ws.Range("E4").Value = 15.03586575318281
$ws.Range("F4").Value = 18.99906936892011
$ws.Range("G4").Value = 19.87665079716189
$ws.Range("H4").Value = 11.93331607525453
$ws.Range("I4").Value = 17.94624457796846
$ws.Range("K4").Value = 10.57913418901649
$ws.Range("N4").Value = 16.82092484775676
$ws.Range("O4").Value = 16.98276762681802
$ws.Range("B5").Value = 6.180019745406446
$ws.Range("D5").Value = 3.569607704589137
$ws.Range("E5").Value = 14.78269186838543
$ws.Range("F5").Value = 18.99024794642647
$ws.Range("G5").Value = 19.84827741308929
$ws.Range("H5").Value = 11.94303160463669
$ws.Range("I5").Value = 17.97193996822383
$ws.Range("K5").Value = 10.4604627238403
$ws.Range("N5").Value = 16.83241111102415
$ws.Range("O5").Value = 16.99199181626541
$ws.Range("B6").Value = 6.174766576110041
$ws.Range("D6").Value = 3.566341235394161
$ws.Range("E6").Value = 14.74029024106473
$ws.Range("F6").Value = 18.98887740894875
$ws.Range("G6").Value = 19.84369574051361
$ws.Range("H6").Value = 11.9446789645601
$ws.Range("I6").Value = 17.97625159527397
$ws.Range("K6").Value = 10.44059217650949
$ws.Range("N6").Value = 16.8343426878227
$ws.Range("O6").Value = 16.99359368659634
$ws.Range("B7").Value = 6.211217537633909
$ws.Range("D7").Value = 3.588897263694295
$ws.Range("E7").Value = 15.03247585096039
$ws.Range("F7").Value = 18.99894408498324
$ws.Range("G7").Value = 19.87625946483428
$ws.Range("H7").Value = 11.93344481543011
$ws.Range("I7").Value = 17.9465881045154
$ws.Range("K7").Value = 10.57754489708062
$ws.Range("N7").Value = 16.82107812759928
$ws.Range("O7").Value = 16.98288731883527
$ws.Range("B8").Value = 6.371476025023521
$ws.Range("D8").Value = 3.68544005722409
$ws.Range("E8").Value = 16.26715251152642
$ws.Range("F8").Value = 19.05702318687232
$ws.Range("G8").Value = 20.03861727335204
$ws.Range("H8").Value = 11.88883469579005
$ws.Range("I8").Value = 17.82220522384409
$ws.Range("K8").Value = 11.15709478124091
$ws.Range("N8").Value = 16.76606964508909
$ws.Range("O8").Value = 16.94587236452099
$ws.Range("B9").Value = 6.680863565973599
$ws.Range("D9").Value = 3.863112945931277
$ws.Range("E9").Value = 18.59153313568305
$ws.Range("F9").Value = 19.21943654535506
$ws.Range("G9").Value = 20.42186547573086
$ws.Range("H9").Value = 11.81925534656191
$ws.Range("I9").Value = 17.60156418249743
$ws.Range("K9").Value = 12.19822312220338
$ws.Range("N9").Value = 16.67081261334683
$ws.Range("O9").Value = 16.91026571754812
$ws.Range("B10").Value = 6.902873369139136
$ws.Range("D10").Value = 3.985826382126144
$ws.Range("E10").Value = 20.24154372726754
$ws.Range("F10").Value = 19.36766088339306
$ws.Range("G10").Value = 20.74059917984116
$ws.Range("H10").Value = 11.77910138542011
$ws.Range("I10").Value = 17.45350837743253
$ws.Range("K10").Value = 12.90054876357363
$ws.Range("N10").Value = 16.60846603872546
$ws.Range("O10").Value = 16.90680186006746
$ws.Range("B11").Value = 7.002265595834681
$ws.Range("D11").Value = 4.039816626884207
$ws.Range("E11").Value = 20.94960852399132
$ws.Range("F11").Value = 19.44118153165799
$ws.Range("G11").Value = 20.89308083470382
$ws.Range("H11").Value = 11.76322593876571
$ws.Range("I11").Value = 17.38917498322668
$ws.Range("K11").Value = 13.20573430007414
$ws.Range("N11").Value = 16.58175097815433
$ws.Range("O11").Value = 16.91017529654811
$ws.Range("B12").Value = 7.03963738568948
$ws.Range("D12").Value = 4.059987404185198
$ws.Range("E12").Value = 21.2116541206746
$ws.Range("F12").Value = 19.4698795146334
$ws.Range("G12").Value = 20.95184453378687
$ws.Range("H12").Value = 11.75755889035392
$ws.Range("I12").Value = 17.36524536581601
$ws.Range("K12").Value = 13.31919376195343
$ws.Range("N12").Value = 16.5718706882139
$ws.Range("O12").Value = 16.91216509008344
$ws.Range("B13").Value = 7.031601071864578
$ws.Range("D13").Value = 4.055655620137943
$ws.Range("E13").Value = 21.15548793967854
$ws.Range("F13").Value = 19.46366107922644
$ws.Range("G13").Value = 20.93914415267689
$ws.Range("H13").Value = 11.7587640498833
$ws.Range("I13").Value = 17.37037985246755
$ws.Range("K13").Value = 13.29485269825108
$ws.Range("N13").Value = 16.57398809428469
$ws.Range("O13").Value = 16.91170486675687
$ws.Range("B14").Value = 7.005345720694035
$ws.Range("D14").Value = 4.041481641352869
$ws.Range("E14").Value = 20.9712889002919
$ws.Range("F14").Value = 19.44352546409418
$ws.Range("G14").Value = 20.89789519469653
$ws.Range("H14").Value = 11.76275279557561
$ws.Range("I14").Value = 17.38719763143257
$ws.Range("K14").Value = 13.21511118097338
$ws.Range("N14").Value = 16.58093339227181
$ws.Range("O14").Value = 16.91032472001121
$ws.Range("B15").Value = 6.989227912740653
$ws.Range("D15").Value = 4.032763665354953
$ws.Range("E15").Value = 20.85767043828776
$ws.Range("F15").Value = 19.43130287767688
$ws.Range("G15").Value = 20.87276047481897
$ws.Range("H15").Value = 11.7652409237108
$ws.Range("I15").Value = 17.39755521429124
$ws.Range("K15").Value = 13.1659913515565
$ws.Range("N15").Value = 16.58521831845988
$ws.Range("O15").Value = 16.90957211777748
$ws.Range("B16").Value = 6.896342537340258
$ws.Range("D16").Value = 3.982260188174246
$ws.Range("E16").Value = 20.19441698050034
$ws.Range("F16").Value = 19.36297716807615
$ws.Range("G16").Value = 20.73078021978504
$ws.Range("H16").Value = 11.78018706566592
$ws.Range("I16").Value = 17.45777327903167
$ws.Range("K16").Value = 12.88031180148098
$ws.Range("N16").Value = 16.61024500903232
$ws.Range("O16").Value = 16.906681058792
$ws.Range("B17").Value = 6.838924224696311
$ws.Range("D17").Value = 3.950800732499625
$ws.Range("E17").Value = 19.77666226829241
$ws.Range("F17").Value = 19.3226092103401
$ws.Range("G17").Value = 20.64555983894022
$ws.Range("H17").Value = 11.78996899899501
$ws.Range("I17").Value = 17.49548669652997
$ws.Range("K17").Value = 12.70135408399864
$ws.Range("N17").Value = 16.62601935228722
$ws.Range("O17").Value = 16.90617575672003
$ws.Range("B18").Value = 6.805750018292735
$ws.Range("D18").Value = 3.932534165331248
$ws.Range("E18").Value = 19.53237541026964
$ws.Range("F18").Value = 19.29996533824849
$ws.Range("G18").Value = 20.59725159592846
$ws.Range("H18").Value = 11.79582029391329
$ws.Range("I18").Value = 17.51746263353587
$ws.Range("K18").Value = 12.59707827840132
$ws.Range("N18").Value = 16.63524737243347
$ws.Range("O18").Value = 16.90635095853337
$ws.Range("B19").Value = 6.794493382413645
$ws.Range("D19").Value = 3.926320228610632
$ws.Range("E19").Value = 19.44897513550695
$ws.Range("F19").Value = 19.29239776351043
$ws.Range("G19").Value = 20.58101851947519
$ws.Range("H19").Value = 11.79784005662809
$ws.Range("I19").Value = 17.52495217499132
$ws.Range("K19").Value = 12.56154310804398
$ws.Range("N19").Value = 16.6383984696442
$ws.Range("O19").Value = 16.90649025664158
$ws.Range("B20").Value = 6.845052158330052
$ws.Range("D20").Value = 3.95416751777897
$ws.Range("E20").Value = 19.82154727386308
$ws.Range("F20").Value = 19.32684709874372
$ws.Range("G20").Value = 20.6545587780471
$ws.Range("H20").Value = 11.78890440399507
$ws.Range("I20").Value = 17.49144264167289
$ws.Range("K20").Value = 12.72054393443358
$ws.Range("N20").Value = 16.6243241058676
$ws.Range("O20").Value = 16.90618132765846
$ws.Range("B21").Value = 7.013065043163806
$ws.Range("D21").Value = 4.045652403684496
$ws.Range("E21").Value = 21.02555745650276
$ws.Range("F21").Value = 19.44941667366128
$ws.Range("G21").Value = 20.90998371869649
$ws.Range("H21").Value = 11.76157184443266
$ws.Range("I21").Value = 17.38224613188737
$ws.Range("K21").Value = 13.23859076332386
$ws.Range("N21").Value = 16.57888698683625
$ws.Range("O21").Value = 16.910710767513
$ws.Range("B22").Value = 7.121306387951414
$ws.Range("D22").Value = 4.103840703654718
$ws.Range("E22").Value = 21.77700548955529
$ws.Range("F22").Value = 19.5345101134963
$ws.Range("G22").Value = 21.08284987932056
$ws.Range("H22").Value = 11.74571760503114
$ws.Range("I22").Value = 17.31339729067499
$ws.Range("K22").Value = 13.56486102099482
$ws.Range("N22").Value = 16.55056721235279
$ws.Range("O22").Value = 16.91782298121953
$ws.Range("B23").Value = 7.063690332311076
$ws.Range("D23").Value = 4.072934431200975
$ws.Range("E23").Value = 21.37917563520094
$ws.Range("F23").Value = 19.48864447309172
$ws.Range("G23").Value = 20.9900639132541
$ws.Range("H23").Value = 11.75399520789627
$ws.Range("I23").Value = 17.34991350595194
$ws.Range("K23").Value = 13.39186491544245
$ws.Range("N23").Value = 16.56555632165725
$ws.Range("O23").Value = 16.91364709745826
$ws.Range("B24").Value = 6.842282226553675
$ws.Range("D24").Value = 3.952645954323473
$ws.Range("E24").Value = 19.80126758043312
$ws.Range("F24").Value = 19.3249293904803
$ws.Range("G24").Value = 20.65048821636332
$ws.Range("H24").Value = 11.78938499867853
$ws.Range("I24").Value = 17.49327004296049
$ws.Range("K24").Value = 12.71187253072132
$ws.Range("N24").Value = 16.62509003099816
$ws.Range("O24").Value = 16.90617735839093
$ws.Range("B25").Value = 6.597935019771597
$ws.Range("D25").Value = 3.816372377128654
$ws.Range("E25").Value = 17.9459643714065
$ws.Range("F25").Value = 19.17036736871204
$ws.Range("G25").Value = 20.31147246279901
$ws.Range("H25").Value = 11.83615631062949
$ws.Range("I25").Value = 17.65877627181464
$ws.Range("K25").Value = 11.92728254203659
$ws.Range("N25").Value = 16.69523706845784
$ws.Range("O25").Value = 16.91592134514634
